$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The final paragraph currently holds both the visible text "...initials (URL)
# etc." and the (invisible) _GoBack bookmark. Split it into a text-only
# paragraph followed by a bookmark-only paragraph, so the new content can be
# inserted between them.
$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n)
$full = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.End - 1)
$splitXml = "<w:p $wNs><w:r><w:t>Check generator works for property names with initials (URL) etc.</w:t></w:r></w:p>" + `
    "<w:p $wNs><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$null = $full.InsertXML($splitXml)

# Re-fetch: the text-only paragraph is now the anchor after which the new
# material goes; the bookmark-only paragraph follows it.
$anchor = $d.Paragraphs.Item($n)

# Simple single-run paragraphs to insert, in order, right after the anchor.
$simpleTexts = @(
    "Check generator works for matchers in the default package.",
    "What happens if you try to generate using a superclass that does not have the expected type parameters and constructors?",
    "Generated class should be public",
    "Generated class should be abstract iff option selected for extensible class",
    "Generated class should be final iff option for extensible class is unselected",
    "Generated class should take type parameters R and T (extends source class) iff it is extensible."
)

$cur = $anchor
foreach ($t in $simpleTexts) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($cur.Index + 1)
    $cur.Range.Text = $t
}

# Paragraph with a lastRenderedPageBreak marker before the text run.
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($cur.Index + 1)
$cur.Range.Text = "PLACEHOLDER"
$target = $d.Range($cur.Range.Start, $cur.Range.End - 1)
$xml = "<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>Generated class should extend selected superclass if specified, passing R and T if extensible and this type and the source class if not.</w:t></w:r></w:p>"
$null = $target.InsertXML($xml)

# Paragraph made up of several runs.
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($cur.Index + 1)
$cur.Range.Text = "PLACEHOLDER"
$target = $d.Range($cur.Range.Start, $cur.Range.End - 1)
$xml = "<w:p $wNs>" + `
    '<w:r><w:t xml:space="preserve">Generated class should extend </w:t></w:r>' + `
    '<w:r><w:t>com.mistraltech.smog.core.CompositePropertyMatcher</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> if</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> no superclass</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> specified, passing R and T if extensible and this type and the source class if not.</w:t></w:r>' + `
    "</w:p>"
$null = $target.InsertXML($xml)

# Finally, append a new empty paragraph after the bookmark-carrying paragraph
# (now the last one in the document).
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($last.Range.End, $last.Range.End)
$null = $endRange.InsertXML("<w:p $wNs/>")
